# Applies the value updates from the scheduled market-price refresh run.
# Each block updates the H/I/J/K/L/M/N (price & profit) columns for one
# leve row, identified by its "Leve Item ID" (column G), across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Item ID 5503)
$ws.Range("H5").Value = 71428840
$ws.Range("J5").Value = 142857520
$ws.Range("L5").Value = 142857520
$ws.Range("N5").Value = -142857750

# Row 40 (Item ID 5505)
$ws.Range("H40").Value = 4249.5
$ws.Range("I40").Value = 4498.5
$ws.Range("J40").Value = 4166.5
$ws.Range("K40").Value = 4498.5
$ws.Range("L40").Value = 4166.5
$ws.Range("M40").Value = -4323.5
$ws.Range("N40").Value = -4516.5

# Row 41 (Item ID 5478)
$ws.Range("H41").Value = 3681
$ws.Range("I41").Value = 3158.6667
$ws.Range("K41").Value = 3158.6667
$ws.Range("M41").Value = -2718.6667

# Row 51 (Item ID 5486)
$ws.Range("H51").Value = 14999.353
$ws.Range("I51").Value = 23428.5
$ws.Range("J51").Value = 9098.950000000001
$ws.Range("K51").Value = 23428.5
$ws.Range("L51").Value = 9098.950000000001
$ws.Range("M51").Value = -22944.5
$ws.Range("N51").Value = -10066.95

# Row 64 (Item ID 5506)
$ws.Range("H64").Value = 7913.4
$ws.Range("I64").Value = 7891.75
$ws.Range("K64").Value = 7891.75
$ws.Range("M64").Value = -7643.75

# Row 67 (Item ID 5506)
$ws.Range("H67").Value = 7913.4
$ws.Range("I67").Value = 7891.75
$ws.Range("K67").Value = 7891.75
$ws.Range("M67").Value = -7033.75

# Row 69 (Item ID 12616)
$ws.Range("H69").Value = 13607.571
$ws.Range("I69").Value = 10490.375
$ws.Range("J69").Value = 15525.846
$ws.Range("K69").Value = 31471.125
$ws.Range("L69").Value = 46577.538
$ws.Range("M69").Value = -30597.125
$ws.Range("N69").Value = -48325.538

# Row 70 (Item ID 12604)
$ws.Range("H70").Value = 35967.668
$ws.Range("J70").Value = 51451.5
$ws.Range("L70").Value = 154354.5
$ws.Range("N70").Value = -154894.5

# Row 72 (Item ID 12616)
$ws.Range("H72").Value = 13607.571
$ws.Range("I72").Value = 10490.375
$ws.Range("J72").Value = 15525.846
$ws.Range("K72").Value = 94413.375
$ws.Range("L72").Value = 139732.614
$ws.Range("M72").Value = -90045.375
$ws.Range("N72").Value = -148468.614

# Row 73 (Item ID 12604)
$ws.Range("H73").Value = 35967.668
$ws.Range("J73").Value = 51451.5
$ws.Range("L73").Value = 154354.5
$ws.Range("N73").Value = -156226.5

# Row 76 (Item ID 12602)
$ws.Range("H76").Value = 1000000000
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# Row 79 (Item ID 12602)
$ws.Range("H79").Value = 1000000000
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# Row 113 (Item ID 27775)
$ws.Range("H113").Value = 6277.2354
$ws.Range("J113").Value = 6491.3335
$ws.Range("L113").Value = 6491.3335
$ws.Range("N113").Value = -12999.3335

# Row 129 (Item ID 36115)
$ws.Range("H129").Value = 1403.5834
$ws.Range("I129").Value = 824.3
$ws.Range("K129").Value = 2472.9
$ws.Range("M129").Value = 2527.1

# Row 134 (Item ID 41997)
$ws.Range("H134").Value = 27228.916
$ws.Range("J134").Value = 27228.916
$ws.Range("L134").Value = 27228.916
$ws.Range("N134").Value = -37368.916

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (Item ID 43999)
$ws.Range("H61").Value = 3870.0715
$ws.Range("I61").Value = 3571.739
$ws.Range("K61").Value = 3571.739
$ws.Range("M61").Value = -3359.739

# Row 63 (Item ID 12528)
$ws.Range("H63").Value = 3233.6365
$ws.Range("J63").Value = 2990
$ws.Range("L63").Value = 2990
$ws.Range("N63").Value = -4362

# Row 66 (Item ID 12528)
$ws.Range("H66").Value = 3233.6365
$ws.Range("J66").Value = 2990
$ws.Range("L66").Value = 14950
$ws.Range("N66").Value = -21814

# Row 136 (Item ID 43999)
$ws.Range("H136").Value = 3870.0715
$ws.Range("I136").Value = 3571.739
$ws.Range("K136").Value = 10715.217
$ws.Range("M136").Value = -8165.217000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Item ID 5092)
$ws.Range("H22").Value = 401.54544
$ws.Range("J22").Value = 150
$ws.Range("L22").Value = 150
$ws.Range("N22").Value = -496

# Row 64 (Item ID 14184)
$ws.Range("H64").Value = 1487.7
$ws.Range("J64").Value = 1889.2
$ws.Range("L64").Value = 1889.2
$ws.Range("N64").Value = -2339.2

# Row 67 (Item ID 14184)
$ws.Range("H67").Value = 1487.7
$ws.Range("J67").Value = 1889.2
$ws.Range("L67").Value = 1889.2
$ws.Range("N67").Value = -3449.2

# Row 94 (Item ID 19939)
$ws.Range("H94").Value = 29604
$ws.Range("I94").Value = 16384.5
$ws.Range("J94").Value = 40179.6
$ws.Range("K94").Value = 16384.5
$ws.Range("L94").Value = 40179.6
$ws.Range("M94").Value = -15933.5
$ws.Range("N94").Value = -41081.6

# Row 134 (Item ID 43998)
$ws.Range("H134").Value = 8479.839
$ws.Range("I134").Value = 8590.115
$ws.Range("K134").Value = 25770.345
$ws.Range("M134").Value = -23235.345

$ws = $wb.Worksheets.Item("CRP")
# Row 94 (Item ID 32934)
$ws.Range("H94").Value = 1431.2273
$ws.Range("I94").Value = 1570.7
$ws.Range("J94").Value = 1315
$ws.Range("K94").Value = 1570.7
$ws.Range("L94").Value = 1315
$ws.Range("M94").Value = -1119.7
$ws.Range("N94").Value = -2217

# Row 122 (Item ID 36196)
$ws.Range("H122").Value = 1549.8422
$ws.Range("I122").Value = 896.3077
$ws.Range("J122").Value = 2965.8333
$ws.Range("K122").Value = 2688.9231
$ws.Range("L122").Value = 8897.499899999999
$ws.Range("M122").Value = -238.9231
$ws.Range("N122").Value = -13797.4999

$ws = $wb.Worksheets.Item("CUL")
# Row 140 (Item ID 44097)
$ws.Range("H140").Value = 7608.3887
$ws.Range("I140").Value = 2919.3076
$ws.Range("J140").Value = 19800
$ws.Range("K140").Value = 8757.9228
$ws.Range("L140").Value = 59400
$ws.Range("M140").Value = -3577.9228
$ws.Range("N140").Value = -69760

$ws = $wb.Worksheets.Item("GSM")
# Row 135 (Item ID 42006)
$ws.Range("H135").Value = 97499
$ws.Range("J135").Value = 97499
$ws.Range("L135").Value = 97499
$ws.Range("N135").Value = -107639

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Item ID 5277)
$ws.Range("H22").Value = 1427.9
$ws.Range("J22").Value = 1469.4783
$ws.Range("L22").Value = 1469.4783
$ws.Range("N22").Value = -2059.4783

# Row 27 (Item ID 5277)
$ws.Range("H27").Value = 1427.9
$ws.Range("J27").Value = 1469.4783
$ws.Range("L27").Value = 1469.4783
$ws.Range("N27").Value = -1683.4783

# Row 40 (Item ID 36248)
$ws.Range("H40").Value = 2590.4119
$ws.Range("I40").Value = 2417.3635
$ws.Range("J40").Value = 2907.6667
$ws.Range("K40").Value = 2417.3635
$ws.Range("L40").Value = 2907.6667
$ws.Range("M40").Value = -2281.3635
$ws.Range("N40").Value = -3179.6667

# Row 68 (Item ID 12563)
$ws.Range("H68").Value = 2511.6
$ws.Range("I68").Value = 2497.4
$ws.Range("J68").Value = 2540
$ws.Range("K68").Value = 2497.4
$ws.Range("L68").Value = 2540
$ws.Range("M68").Value = -1748.4
$ws.Range("N68").Value = -4038

# Row 71 (Item ID 12563)
$ws.Range("H71").Value = 2511.6
$ws.Range("I71").Value = 2497.4
$ws.Range("J71").Value = 2540
$ws.Range("K71").Value = 12487
$ws.Range("L71").Value = 12700
$ws.Range("M71").Value = -8743
$ws.Range("N71").Value = -20188

$ws = $wb.Worksheets.Item("WVR")
# Row 28 (Item ID 3053)
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# Row 96 (Item ID 19977)
$ws.Range("H96").Value = 1245.8
$ws.Range("I96").Value = 1307.25
$ws.Range("K96").Value = 1307.25
$ws.Range("M96").Value = 65.75
